$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 49380
$ws.Range("B15").Value = 121
$ws.Range("C15").Value = 'Cluj-Turda'
$ws.Range("D15").Value = 'Interes Serviciu'
$ws.Range("B16").Value = 421
$ws.Range("C16").Value = 'Cluj-Satu-Mare'
$ws.Range("D16").Value = 'Interes Serviciu'
$ws.Range("B19").Value = 85
$ws.Range("C19").Value = 'Cluj-Apahida'
$ws.Range("B21").Value = 101
$ws.Range("C21").Value = 'Cluj-Dej'
$ws.Range("B22").Value = 47
$ws.Range("C22").Value = 'Cluj-Cluj'
$ws.Range("B23").Value = 356
$ws.Range("C23").Value = 'Cluj-Baia-Mare'
$ws.Range("D23").Value = 'Interes Serviciu'
$ws.Range("B27").Value = 92
$ws.Range("C27").Value = 'Cluj-Bontida'
$ws.Range("B33").Value = 101
$ws.Range("C33").Value = 'Cluj-Dej'
$ws.Range("B34").Value = 156
$ws.Range("C34").Value = 'Cluj-Zalau'
$ws.Range("B35").Value = 257
$ws.Range("C35").Value = 'Cluj-Bistrita'
$ws.Range("D35").Value = 'Interes Serviciu'
$ws.Range("B36").Value = 30
$ws.Range("C36").Value = 'Acasa-Birou'
$ws.Range("D36").Value = ' '
$ws.Range("B37").Value = 257
$ws.Range("C37").Value = 'Cluj-Bistrita'
$ws.Range("D37").Value = 'Interes Serviciu'
$ws.Range("B40").Value = 156
$ws.Range("C40").Value = 'Cluj-Zalau'
$ws.Range("B42").Value = 30
$ws.Range("C42").Value = 'Acasa-Birou'
$ws.Range("D42").Value = ' '
$ws.Range("B44").Value = 2895
$ws.Range("B45").Value = 52275
